# "update doc and example" — refresh the MaxSAT example worksheets/charts with
# corrected sample data, tidy up the saved cursor position on each example
# sheet, and set up the page for printing on the first example sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Example 1 — swap a few mixed-up parking id values in column A
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Example 1")

$ws1.Range("A7").Value = 7
$ws1.Range("A8").Value = 5

$ws1.Range("A11").Value = 7
$ws1.Range("A12").Value = 4

$ws1.Range("A16").Value = 5

$ws1.Range("A19").Value = 7
$ws1.Range("A20").Value = 4

# Page setup for printing (Page Setup dialog -> Paper size A4, Portrait)
$ws1.PageSetup.PaperSize = 9
$ws1.PageSetup.Orientation = 1

# Reset the saved selection back to the top-left of the sheet.
$ws1.Activate() | Out-Null
$ws1.Range("A1").Select() | Out-Null

# ---------------------------------------------------------------------------
# Example 2 — correct a couple of sampled values, move the saved cursor
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Example 2")

$ws2.Range("A12").Value = 34
$ws2.Range("B14").Value = 15
$ws2.Range("A20").Value = 34

$ws2.Activate() | Out-Null
$ws2.Range("L30").Select() | Out-Null

# ---------------------------------------------------------------------------
# Example 3 — correct a couple of sampled values, move the saved cursor
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Example 3")

$ws3.Range("B14").Value = 17
$ws3.Range("A18").Value = 93
$ws3.Range("A19").Value = 94
$ws3.Range("B19").Value = 14
$ws3.Range("A20").Value = 92

$ws3.Activate() | Out-Null
$ws3.Range("M24").Select() | Out-Null

# ---------------------------------------------------------------------------
# Example 4 — correct a couple of sampled values, move the saved cursor
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Example 4")

$ws4.Range("A8").Value = 69
$ws4.Range("A20").Value = 27

$ws4.Activate() | Out-Null
$ws4.Range("Q23").Select() | Out-Null

# Leave the first example sheet as the active tab, matching the saved file.
$ws1.Activate() | Out-Null
